$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.941.69"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.046.15"

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.00"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.05"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("E12").Value = "  +3.62%  "

$ws.Range("D13").Value = "2.340.76"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("E14").Value = "  +5.97%  "

$ws.Range("E15").Value = "  -4.39%  "

$ws.Range("D16").Value = "2.040.88"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "36.892.29"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.24"
$ws.Range("E18").Value = "  +13.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.96"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.31"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.73"
$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +8.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.45"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.04"
$ws.Range("E27").Value = "  -1.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.75"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("E30").Value = "  +4.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.65"
$ws.Range("E31").Value = "  +0.82%  "

$ws.Range("E32").Value = "  -2.59%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E35").Value = "  +3.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.22"
$ws.Range("E36").Value = "  -3.16%  "

$ws.Range("E37").Value = "  -1.88%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.33"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("E39").Value = "  -5.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +13.53%  "

$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.90"
$ws.Range("E41").Value = "  +24.44%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0220"
$ws.Range("E42").Value = "  -2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.17"
$ws.Range("E43").Value = "  -5.88%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.32"
$ws.Range("E44").Value = "  -2.08%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -3.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"

$ws.Range("D47").Value = "1.277.97"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  -3.29%  "

$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.43"
$ws.Range("E51").Value = "  -22.40%  "
